# The two species records that occupy rows 2 and 3 were swapped: the
# "Tretåig hackspett" (Picoides tridactylus) record moves from row 2 to
# row 3, and the "Trådticka" (Climacocystis borealis) record moves from
# row 3 to row 2. Columns O onward already hold identical values on both
# rows, so only columns A-N plus AF (which carry the species-specific
# data) need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes the "Trådticka" / Climacocystis borealis record ---
$ws.Range("A2").Value = 111813166
$ws.Range("B2").Value = 90087
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 3298
$ws.Range("F2").Value = "Trådticka"
$ws.Range("G2").Value = "Climacocystis borealis"
$ws.Range("H2").Value = "(Fr.) Kotl. & Pouzar"

# row2 gains the (empty) "Enhet" cell (J2) ...
$ws.Range("J2").Value = ""
# ... and loses the "Kön" (L2) / "Aktivitet" (M2) cells it used to carry
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# row2 gains the (empty) "Bestämningsmetod" cell (AF2)
$ws.Range("AF2").Value = ""

# --- Row 3 becomes the "Tretåig hackspett" / Picoides tridactylus record ---
$ws.Range("A3").Value = 111813153
$ws.Range("B3").Value = 56398
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

# row3 loses the (empty) "Enhet" cell (J3) ...
$ws.Range("J3").ClearContents()
# ... and gains the "Kön" (L3, empty) / "Aktivitet" (M3) cells
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "färska spår"

# row3 loses the (empty) "Bestämningsmetod" cell (AF3)
$ws.Range("AF3").ClearContents()
